{"js": "// Replace the two-digit/one-digit division expressions in the worksheet\n// table cells with their updated values, as produced by the authoring\n// edit captured in the commit. Each left-hand expression occurs exactly\n// once in the document, so a targeted search/replace per pair is safe.\nconst replacements = [\n  [\"72\u00f74=\", \"90\u00f76=\"],\n  [\"15\u00f72=\", \"86\u00f77=\"],\n  [\"95\u00f76=\", \"67\u00f78=\"],\n  [\"12\u00f73=\", \"14\u00f78=\"],\n  [\"84\u00f75=\", \"43\u00f74=\"],\n  [\"40\u00f79=\", \"64\u00f79=\"],\n  [\"12\u00f76=\", \"77\u00f76=\"],\n  [\"13\u00f72=\", \"53\u00f79=\"],\n  [\"58\u00f75=\", \"29\u00f79=\"],\n  [\"44\u00f74=\", \"37\u00f76=\"],\n  [\"24\u00f77=\", \"67\u00f76=\"],\n  [\"62\u00f77=\", \"36\u00f79=\"],\n  [\"27\u00f79=\", \"21\u00f76=\"],\n  [\"30\u00f78=\", \"49\u00f72=\"],\n  [\"15\u00f76=\", \"58\u00f74=\"],\n  [\"59\u00f76=\", \"77\u00f73=\"],\n  [\"12\u00f77=\", \"53\u00f73=\"],\n  [\"84\u00f77=\", \"37\u00f78=\"],\n  [\"52\u00f78=\", \"57\u00f73=\"],\n  [\"28\u00f74=\", \"75\u00f76=\"],\n  [\"57\u00f75=\", \"74\u00f73=\"],\n  [\"67\u00f73=\", \"31\u00f75=\"],\n  [\"76\u00f77=\", \"37\u00f76=\"],\n  [\"31\u00f73=\", \"32\u00f75=\"],\n  [\"61\u00f77=\", \"51\u00f75=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit/one-digit division expressions in the worksheet\n# table to the regenerated values, as captured by the authoring commit.\n# Each \"Find\" expression occurs exactly once in the document, so a\n# targeted Find/Replace per pair is safe and avoids touching anything\n# else (date heading, cell formatting, etc.).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"72\u00f74=\"; Replace = \"90\u00f76=\" },\n    @{ Find = \"15\u00f72=\"; Replace = \"86\u00f77=\" },\n    @{ Find = \"95\u00f76=\"; Replace = \"67\u00f78=\" },\n    @{ Find = \"12\u00f73=\"; Replace = \"14\u00f78=\" },\n    @{ Find = \"84\u00f75=\"; Replace = \"43\u00f74=\" },\n    @{ Find = \"40\u00f79=\"; Replace = \"64\u00f79=\" },\n    @{ Find = \"12\u00f76=\"; Replace = \"77\u00f76=\" },\n    @{ Find = \"13\u00f72=\"; Replace = \"53\u00f79=\" },\n    @{ Find = \"58\u00f75=\"; Replace = \"29\u00f79=\" },\n    @{ Find = \"44\u00f74=\"; Replace = \"37\u00f76=\" },\n    @{ Find = \"24\u00f77=\"; Replace = \"67\u00f76=\" },\n    @{ Find = \"62\u00f77=\"; Replace = \"36\u00f79=\" },\n    @{ Find = \"27\u00f79=\"; Replace = \"21\u00f76=\" },\n    @{ Find = \"30\u00f78=\"; Replace = \"49\u00f72=\" },\n    @{ Find = \"15\u00f76=\"; Replace = \"58\u00f74=\" },\n    @{ Find = \"59\u00f76=\"; Replace = \"77\u00f73=\" },\n    @{ Find = \"12\u00f77=\"; Replace = \"53\u00f73=\" },\n    @{ Find = \"84\u00f77=\"; Replace = \"37\u00f78=\" },\n    @{ Find = \"52\u00f78=\"; Replace = \"57\u00f73=\" },\n    @{ Find = \"28\u00f74=\"; Replace = \"75\u00f76=\" },\n    @{ Find = \"57\u00f75=\"; Replace = \"74\u00f73=\" },\n    @{ Find = \"67\u00f73=\"; Replace = \"31\u00f75=\" },\n    @{ Find = \"76\u00f77=\"; Replace = \"37\u00f76=\" },\n    @{ Find = \"31\u00f73=\"; Replace = \"32\u00f75=\" },\n    @{ Find = \"61\u00f77=\"; Replace = \"51\u00f75=\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Find, $false, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, 2)\n}\n"}
